# Update the worksheet data: new empadronador totals, reordered names, and
# remove the last row (row 12) since the data now only spans A1:B11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2 through 11 (name, total_registros)
$data = @(
    @("CARRION LAZARO MICHAEL LUIS", 152),
    @("ARRUNATEGUI ESPINOZA JOVANNY", 148),
    @("NIMA CARMEN KAREN DEL MILAGRO", 134),
    @("PAZ ANASTACIO JUANITA ROSA", 133),
    @("ESPINOZA VALDIVIEZO JUNIOR RICARDO", 121),
    @("ALZAMORA CHERRES SIRLEY YASMIN", 114),
    @("PULACHE LAZO VILMA YOHANA", 112),
    @("DOMINGUEZ CUEVA MERLING DEL JESUS YOLINDA", 109),
    @("NAVARRO JUAREZ LIDIA", 108),
    @("LILIAN ROXANA VEGA GARCÍA", 104)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $row = $row + 1
}

# Remove the previous row 12 contents (NIMA CRUZ ANA GRACIELA | 1); the data
# range is now A1:B11, so delete the whole row 12 to shrink the sheet.
$ws.Rows.Item(12).Delete()
